# "Generate Report for Handback"
#
# The localization-status report gets a handback pass: the status text
# flips from "Ready for handoff" to "Handed back: in sync with en-US"
# (wherever that status string appears), the per-language "Latest
# Handback DateTime" cells get real timestamps instead of the
# 0001-01-01 placeholder, and each language sheet grows a "Latest
# Target File" / "Latest Handback File" hyperlink pair (columns F/G)
# mirroring the existing Source/Handoff-file hyperlinks in columns A/D.

$wb = $excel.ActiveWorkbook

$Overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# 1. Status flips from "Ready for handoff" to "Handed back: in sync
#    with en-US" everywhere it shows up (Overview summary columns plus
#    each language sheet's Status column).
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$Overview.Range("B2").Value = $newStatus
$Overview.Range("C2").Value = $newStatus
$Overview.Range("B3").Value = $newStatus
$Overview.Range("C3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Latest Handback DateTime (column H) now carries the real handback
#    timestamp for each language instead of the 0001-01-01 placeholder.
# ---------------------------------------------------------------------
$zhcn.Range("H2").Value = "2016-03-24 14:42:42"
$zhcn.Range("H3").Value = "2016-03-24 14:42:42"

$dede.Range("H2").Value = "2016-03-24 14:42:51"
$dede.Range("H3").Value = "2016-03-24 14:42:51"

# ---------------------------------------------------------------------
# 3. New "Latest Target File" (F) / "Latest Handback File" (G) columns:
#    same source-file / translated-file hyperlink pair already shown in
#    A (Source File Name) and D (Latest Handoff File), now duplicated
#    to reflect the handback artifacts.
# ---------------------------------------------------------------------
$zhcnXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$dedeXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$aMdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/15df40e4cb62275e0c02ad80e091be7e2176fc4a/e2e/a.md"
$zhcnXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8103497a9b740b6839ea997910ec3d9fd806790e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$dedeXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/57a6004a0e9e6f2b581bc96e7b263fca8d016c47/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# zh-cn sheet: rows 2 and 3
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), $zhcnXlfUrl, [Type]::Missing, [Type]::Missing, $zhcnXlf) | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), $zhcnXlfUrl, [Type]::Missing, [Type]::Missing, $zhcnXlf) | Out-Null

# de-de sheet: rows 2 and 3
$dede.Hyperlinks.Add($dede.Range("F2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G2"), $dedeXlfUrl, [Type]::Missing, [Type]::Missing, $dedeXlf) | Out-Null
$dede.Hyperlinks.Add($dede.Range("F3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G3"), $dedeXlfUrl, [Type]::Missing, [Type]::Missing, $dedeXlf) | Out-Null

"Generated handback report"
